$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Schedule" updates (rows 2-4)
# ----------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

# Row 2
$schedule.Range("A2").Value = 46053.25
$schedule.Range("C2").Value = 13
$schedule.Range("D2").Value = 49.14
$schedule.Range("E2").Value = 1047.58868175
$schedule.Range("F2").Value = 21.31845099206349

# Row 3
$schedule.Range("A3").Value = 46053.875
$schedule.Range("B3").Value = 46054.08333333334
$schedule.Range("C3").Value = 5
$schedule.Range("D3").Value = 18.9
$schedule.Range("E3").Value = 781.9951815000001
$schedule.Range("F3").Value = 41.37540642857144

# Row 4
$schedule.Range("A4").Value = 46054.25
$schedule.Range("C4").Value = 10
$schedule.Range("D4").Value = 37.8
$schedule.Range("E4").Value = 862.4053132500001
$schedule.Range("F4").Value = 22.81495537698413

# ----------------------------------------------------------------------
# Sheet "Detailed" updates
# ----------------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("E14").Value = "ON"

$detailed.Range("B37").Value = 73.55727

$detailed.Range("B38").Value = 103.98809

$detailed.Range("B39").Value = 105
$detailed.Range("C39").Value = "historical"

$detailed.Range("B40").Value = 147.52
$detailed.Range("C40").Value = "historical"

$detailed.Range("B41").Value = 146.76862
$detailed.Range("C41").Value = "historical"

$detailed.Range("B42").Value = 147.52
$detailed.Range("C42").Value = "historical"

$detailed.Range("B43").Value = 115.25453
$detailed.Range("C43").Value = "historical"
$detailed.Range("E43").Value = "OFF"

$detailed.Range("B44").Value = 101.25
$detailed.Range("C44").Value = "historical"

$detailed.Range("B45").Value = 105
$detailed.Range("C45").Value = "historical"

$detailed.Range("B46").Value = 98.23987
$detailed.Range("C46").Value = "historical"

$detailed.Range("B47").Value = 78
$detailed.Range("C47").Value = "historical"

$detailed.Range("B48").Value = 70.38733999999999
$detailed.Range("C48").Value = "historical"

$detailed.Range("B49").Value = 63.44156

$detailed.Range("B50").Value = 57.31

$detailed.Range("B51").Value = 78
$detailed.Range("E51").Value = "ON"

$detailed.Range("B52").Value = 76.12041000000001
$detailed.Range("E52").Value = "ON"

$detailed.Range("B53").Value = 74.29716000000001
$detailed.Range("E53").Value = "ON"

$detailed.Range("B54").Value = 73.20010000000001

$detailed.Range("B55").Value = 73.20010000000001

$detailed.Range("B56").Value = 72.51924

$detailed.Range("B57").Value = 63.1318

$detailed.Range("B58").Value = 61.57904

$detailed.Range("B59").Value = 78.73907
$detailed.Range("E59").Value = "OFF"

$detailed.Range("B60").Value = 75.28270000000001
$detailed.Range("E60").Value = "OFF"

$detailed.Range("B61").Value = 76.99852
$detailed.Range("E61").Value = "OFF"

$detailed.Range("B62").Value = 73.20001999999999

$detailed.Range("B63").Value = 64.97036

$detailed.Range("B67").Value = 35.88

$detailed.Range("B68").Value = 35.88

$detailed.Range("B73").Value = 35.88

$detailed.Range("B78").Value = 36.06071

$detailed.Range("B79").Value = 57.06003

$detailed.Range("B80").Value = 59.21238

$detailed.Range("B81").Value = 68.71477

$detailed.Range("B82").Value = 84.79000000000001

$detailed.Range("B83").Value = 78

$detailed.Range("B84").Value = 92.14239000000001

$detailed.Range("B85").Value = 89.95526

$detailed.Range("B86").Value = 80.93996

$detailed.Range("B87").Value = 100.3

$detailed.Range("B88").Value = 147.52

$detailed.Range("B89").Value = 135.04632

$detailed.Range("B91").Value = 82.70392

$detailed.Range("B92").Value = 90.98444000000001

$detailed.Range("B93").Value = 77.03009

$detailed.Range("B94").Value = 57.31

$detailed.Range("B95").Value = 63.73481

$detailed.Range("B96").Value = 57.31

$detailed.Range("B97").Value = 64.81568
